$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target state for rows 2-9 (name, latitude, longitude)
$data = @(
    @("wisdom",      27.1462402,          36.4825152),
    @("gospel",      33.80775722088888,   31.63713273881481),
    @("lyric",       30.95328428020835,   35.43131750034723),
    @("historical",  31.9524206297491,    36.09865291409798),
    @("letter",      34.50393560222223,   29.1487569437037),
    @("prophecy",    32.04909800329218,   35.39011386927299),
    @("apocalyptic", 32.24773232303031,   35.54896866929293),
    @("law",         30.5217654,          34.112561)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}
